# Added the excel data reader code
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$data = @(
    @("Object Name", "Type", "Value"),
    @("table1", "list", "nim`$`$sharma`$`$gaurav`$`$kumar"),
    @("obj2", "label", "kk"),
    @("obj3", "list", "nim`$`$sharma`$`$gaurav`$`$kumar"),
    @("obj4", "label", "tt")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Range("E11").Select() | Out-Null
